# Add support for mono camera
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Reduce the default "bad frame fraction" from 0.2 to 0.05
$ws.Range("B11").Value = 0.05

# The "黑白相机配滤镜必须选True" hint row: change the True/False selector in
# E16 to False, reusing the exact formatting/type used by identical cells
# nearby (so it stays a shared-string "False" rather than becoming a
# native boolean).
$ws.Range("E13").Copy() | Out-Null
$ws.Range("E16").PasteSpecial(-4163) | Out-Null
$ws.Application.CutCopyMode = $false

# Add the new helper note next to it, explaining that mono cameras with a
# filter wheel must have this set to True.
$ws.Range("G16").Value = "黑白相机配滤镜必须选True"

# Match the author's final selection location.
$ws.Range("G17").Select() | Out-Null

Write-Host "done"
